# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# The "Periodo Mora" column (E) for the worker JOHANA PATRICIA RODRIGUEZ SILVA
# (rows 16-27) is refreshed: the previous list of pending periods (1706..1805,
# ascending) is replaced by the new list, now stored in descending order
# (1805..1706) - i.e. the old account-statement periods are dropped and the
# new ones are appended, updating the underlying database/table order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("1805", "1804", "1803", "1802", "1801", "1712", "1711", "1710", "1709", "1708", "1707", "1706")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value2 = $periodos[$i]
}
